$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.498.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.249.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.248.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.499'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.66%  '
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.784.96'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.638.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.248.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '506.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.750'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("E25").Value = '  +2.88%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.129'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +46.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.92'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.97'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("E35").Value = '  -3.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("E37").Value = '  +20.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0785'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +15.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '493.37'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.88%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.292'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.985.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.32%  '
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.53%  '
